$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data point was recorded for "Albahaca" (Agrícola del Norte S.A. de Arica).
# It belongs right after the header/most-recent block, at row 19, pushing every
# existing row from 19 downwards by one (old row 19 -> new row 20, ..., old row 64 -> new row 65).

$ws.Rows.Item(19).Insert()

# Fill the newly inserted row 19 with the new record. Most attributes repeat the
# values that used to sit in the old row 19 (now shifted to row 20); only the
# date (Fecha) and volume (Volumen) differ for this new entry.
$ws.Range("A19").Value = 1
$ws.Range("B19").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C19").Value = 'Arica y Parinacota'
$ws.Range("D19").Value = 45099
$ws.Range("E19").Value = 15
$ws.Range("F19").Value = 100112052
$ws.Range("G19").Value = 'Albahaca'
$ws.Range("H19").Value = 'Sin especificar'
$ws.Range("I19").Value = 'Primera'
$ws.Range("J19").Value = 300
$ws.Range("K19").Value = 900
$ws.Range("L19").Value = 1000
$ws.Range("M19").Value = 950
$ws.Range("N19").Value = '$/paquete'
$ws.Range("O19").Value = 'Región de Arica y Parinacota'
$ws.Range("P19").Value = 950
$ws.Range("Q19").Value = 1
$ws.Range("R19").Value = 'Hortaliza'
